# Update the "Percent Change" column (E) values on Sheet1 for rows 2-13.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

$values = @{
    2  = 0.003241491085899328
    3  = 0.006236507555768878
    4  = 0.004961020552799456
    5  = 0.002013318878736348
    6  = -0.01795580110497241
    7  = -0.000851996105160624
    8  = 0.007393715341959206
    9  = 0.01048951048951063
    10 = 0.0005783689994216523
    11 = 0.003755006675567385
    12 = -0.006180694419196442
    13 = 0.001878716753191512
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}

$ws.Protect()
